$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.958.99"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.55%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.916.22"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.50%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "608.15"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.97"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +4.70%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.917.97"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.56%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.99%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.40"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.64%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.11%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +4.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.37"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.576.15"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.48%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.911.99"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.79%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.945.29"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.70"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +9.54%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.62"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.25%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.97%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.18"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "493.44"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.80%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +3.75%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +5.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.72"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.30"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.31"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.14"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.88%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.16%  "
$ws.Range("B31").Value = "WrappedeETH"
$ws.Range("C31").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.068.78"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.53%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.43"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.76%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "32.14"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.23%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.879.70"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.90%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.76%  "
$ws.Range("B37").Value = "Mantle"
$ws.Range("C37").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.05"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.60%  "
$ws.Range("B38").Value = "Filecoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.13"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +4.26%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.30"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +11.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.33%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +7.31%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "436.79"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "48.22"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.57%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.83%  "
$ws.Range("B49").Value = "Arweave"
$ws.Range("C49").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "40.80"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +5.64%  "
$ws.Range("B50").Value = "FLOKI"
$ws.Range("C50").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000275"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +20.89%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "143.10"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.24%  "
